# Sync attendance_reports, modules_schedules, and assets from main repo - 2026-01-17 09:15:08
#
# This script reproduces, via Excel COM interop, the same edits that the
# upstream sync applied to the "Session Analysis Results" worksheet:
#   1. Refresh the top Class-Statistics counters (Missing / Pending sessions).
#   2. Re-order the "<email>, System" attendance-taker label to
#      "System, <email>" everywhere it appears (System re-confirmed the
#      session after the instructor, so System is now listed first).
#   3. Swap the Missing/Pending split for the B1 group-statistics rows that
#      were re-classified (rows 15-20).
#   4. Flip the six still-pending B1 sessions (17/01/2026) from "Pending"
#      (yellow highlight) to "Not Recorded" (pink highlight) now that the
#      session date has passed without any attendance being taken.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- 1. Class Statistics (top-right mini table) ---------------------------
$ws.Range("L7").Value = 39   # Missing Sessions
$ws.Range("L8").Value = 30   # Pending Sessions

# --- 2. "<email>, System" -> "System, <email>" reorder ---------------------
$attendanceRows = @(8,9,10,12,14,15,17,18,34,35,36,38,40,41,43,44,60,61,62,64,66,67,69,70,86,87,88,90,92,93,95,96,112,113,114,116,118,119,121,122,138,139,140,142,144,145,147,148,164,167,170,174,191,194,197,201,218,221,224,228,245,248,251,255,272,275,278,282,299,302,305,309)
foreach ($r in $attendanceRows) {
    $ws.Range("G$r").Value = "System, dnasr281@gmail.com"
}

# --- 3. Group Statistics Missing/Pending split (rows 15-20) ----------------
$ws.Range("P15").Value = 3
$ws.Range("Q15").Value = 2

$ws.Range("P16").Value = 2
$ws.Range("Q16").Value = 2

$ws.Range("P17").Value = 2
$ws.Range("Q17").Value = 2

$ws.Range("P18").Value = 2
$ws.Range("Q18").Value = 2

$ws.Range("P19").Value = 2
$ws.Range("Q19").Value = 2

$ws.Range("P20").Value = 3
$ws.Range("Q20").Value = 2

# --- 4. "Pending" (yellow) -> "Not Recorded" (pink) for the 17/01/2026 rows -
# Row 3 already carries the "Not Recorded" (pink) formatting we need, so
# copy its look-and-feel into each target row (columns A-I only) and then
# fix up the status label text.
$formatSource = $ws.Range("A3:I3")
$pendingRows = @(25,51,77,103,129,155)
foreach ($r in $pendingRows) {
    $formatSource.Copy()
    $ws.Range("A" + $r + ":I" + $r).PasteSpecial(-4122)
    $ws.Range("I$r").Value = "Not Recorded"
}

$excel.CutCopyMode = $false
